$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add notes (column F) to a handful of papers, referenced by the row
# they currently occupy (before the re-sort performed below).
$ws.Range("F2").Value = "Suspicious results. Model is decagon, extended with relation attention module."
$ws.Range("F3").Value = "Enzyme and transporter data, also removed ~100 drugs from decagon data"
$ws.Range("F4").Value = "800 epochs, uses enzyme and drug transporter data, "
$ws.Range("F5").Value = "SAME DATA/AUTHORS AS MS-ADR"
$ws.Range("F7").Value = "Chemical substructure data"
$ws.Range("F26").Value = "120 epochs, uses atomic structure data in 'graph of graphs' format"

# Re-sort the table by column C (AUPRC) descending instead of column D.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C1:C35"), $null, 2, $null, 0)
$ws.Sort.SetRange($ws.Range("A2:F35"))
$ws.Sort.Apply()

# Leave the active selection where the editor last left it.
$ws.Range("F8").Select() | Out-Null
